$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Toolhead section: microswitch quantity 1 -> 2 (row 40, column C) ---
$ws.Range("C40").Value = 2

# --- Dock section updates (rows 54-68 before edit) ---

# Row 55: hardware changed from "m3 x 8mm bhcs" to "m3 x 8mm shcs"
# (also the left-border style on column A switches from the "Rear" style to the
#  "Front" style, matching style index 9 used elsewhere in the hardware column)
$ws.Range("B55").Value = "m3 x 8mm shcs"
$ws.Range("A56:A56").Copy()
$ws.Range("A55").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A4").Copy()
$ws.Range("A55").PasteSpecial(-4122)  # xlPasteFormats (style 9, left border no extra alignment)

# Insert a new hardware row for "m3 x 12mm bhcs" above what is currently row 57
$ws.Rows.Item(57).Insert()
$ws.Range("G57").Clear()
$ws.Range("A58:C58").Copy()
$ws.Range("A57:C57").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E58:F58").Copy()
$ws.Range("E57:F57").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B57").Value = "m3 x 12mm bhcs"
$ws.Range("C57").Value = 1

# Existing "m3 brass insert" row (now row 58): quantity 3 -> 1
$ws.Range("C58").Value = 1

# Existing "m3 hex nut" row (now row 59): quantity 3 -> 2
$ws.Range("C59").Value = 2

# Insert a new hardware row for "m3 washer" above what is currently row 60
$ws.Rows.Item(60).Insert()
$ws.Range("G60").Clear()
$ws.Range("A61:C61").Copy()
$ws.Range("A60:C60").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E61:F61").Copy()
$ws.Range("E60:F60").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B60").Value = "m3 washer"
$ws.Range("C60").Value = 1

Write-Output "Edit applied"
